# Apply "#5: cash & deposit done" change to the 存款 (deposit) sheet.
#
# The deposit sheet gains 8 new trailing columns (G:M) that carry the
# normalized/flattened metadata the scraper now emits for every row
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index), plus a proper "total" header label in F1 (which
# used to just duplicate the first data row's number). Row 4's F column
# also gets corrected from the stray text value "145.TO0" to the plain
# number 145.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Header row (row 1) ----------------------------------------------
# The whole header row switches from ad-hoc Chinese sample values
# (leftover from row 2 having been copy/pasted into row 1) to the
# English field-name labels used by the flattened export schema. F1
# used to hold a stray numeric duplicate (46845); it becomes the
# "total" header label. G1:M1 are brand-new header cells, styled like
# B1:E1.
$headers = @{
    "B1" = "bank"
    "C1" = "deposit_type"
    "D1" = "currency"
    "E1" = "owner"
    "F1" = "total"
    "G1" = "property_category"
    "H1" = "category"
    "I1" = "date"
    "J1" = "legislator_name"
    "K1" = "legislator_id"
    "L1" = "source_file"
    "M1" = "index"
}
foreach ($addr in $headers.Keys) {
    $ws.Range("B1").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).Value = $headers[$addr]
}

# ---- Fix row 4's total: was the text "145.TO0", should be number 145 --
$ws.Range("F4").Value = 145

# ---- New trailing columns for every data row (2-10) --------------------
$indexValues = @{
    2 = 42
    3 = 43
    4 = 46
    5 = 47
    6 = 48
    7 = 49
    8 = 50
    9 = 51
    10 = 52
}

for ($r = 2; $r -le 10; $r++) {
    $ws.Range("G$r").Value = "deposit"
    $ws.Range("H$r").Value = "normal"

    # "date" must stay literal text ("2011-12-31"), not get auto-parsed
    # into a date serial by COM's type inference. Force text format,
    # assign, then strip the now-unneeded text format back off by
    # re-pasting the (style-less) formats from the cell we just wrote,
    # so the cell ends up with the same "no explicit style" look as its
    # siblings.
    $ws.Range("I$r").NumberFormat = "@"
    $ws.Range("I$r").Value = "2011-12-31"
    $ws.Range("G$r").Copy() | Out-Null
    $ws.Range("I$r").PasteSpecial(-4122) | Out-Null

    $ws.Range("J$r").Value = "劉建國"
    $ws.Range("K$r").Value = 1723
    $ws.Range("L$r").Value = "tmp9aab1"
    $ws.Range("M$r").Value = $indexValues[$r]
}
